$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns to Text format so numeric-looking strings
# (e.g. "933,219.00", "0.66") are stored as literal text, matching
# the source data (General style, t="s" cells with no number formatting).
$textCols = @("A","B","C","E","F","G","K","L")
foreach ($col in $textCols) {
    $ws.Range(($col + "1:" + $col + "18")).NumberFormat = "@"
}

$ws.Range("A1").Value = "Collector"
$ws.Range("B1").Value = "Team"
$ws.Range("C1").Value = "Cycle"
$ws.Range("D1").Value = "Repayment_collections"
$ws.Range("E1").Value = "Repayment_amount"
$ws.Range("F1").Value = "Pending Amount"
$ws.Range("G1").Value = "Pending Amount Recovery"
$ws.Range("H1").Value = "Talk_time"
$ws.Range("I1").Value = "New_collections"
$ws.Range("J1").Value = "Repayment_new_collections"
$ws.Range("K1").Value = "New_collection_amount_rate"
$ws.Range("L1").Value = "New_collection_count_rate"
$ws.Range("A2").Value = "Ridhoi Berkat Zebua"
$ws.Range("B2").Value = "Hansyah_S2l"
$ws.Range("C2").Value = "S2"
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "2,991,200.00"
$ws.Range("F2").Value = "158,633,067.00"
$ws.Range("G2").Value = "1.89"
$ws.Range("H2").Value = 1.5569999999999999
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "0.00"
$ws.Range("L2").Value = "0.00"
$ws.Range("A3").Value = "Adistira Winditya P"
$ws.Range("B3").Value = "Hansyah_S2l"
$ws.Range("C3").Value = "S2"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "933,219.00"
$ws.Range("F3").Value = "140,905,074.00"
$ws.Range("G3").Value = "0.66"
$ws.Range("H3").Value = 822
$ws.Range("I3").Value = 15
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "1.55"
$ws.Range("L3").Value = "6.67"
$ws.Range("A4").Value = "Yandi Nugraha"
$ws.Range("B4").Value = "Hansyah_S2l"
$ws.Range("C4").Value = "S2"
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = "2,285,432.00"
$ws.Range("F4").Value = "112,590,060.00"
$ws.Range("G4").Value = "2.03"
$ws.Range("H4").Value = 688
$ws.Range("I4").Value = 14
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = "10.29"
$ws.Range("L4").Value = "7.14"
$ws.Range("A5").Value = "Wasti Feronika Sihombing"
$ws.Range("B5").Value = "Hansyah_S2l"
$ws.Range("C5").Value = "S2"
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = "655,679.00"
$ws.Range("F5").Value = "145,384,324.00"
$ws.Range("G5").Value = "0.45"
$ws.Range("H5").Value = 721
$ws.Range("I5").Value = 14
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = "4.11"
$ws.Range("L5").Value = "0.00"
$ws.Range("A6").Value = "Sucika Wardani"
$ws.Range("B6").Value = "Hansyah_S2l"
$ws.Range("C6").Value = "S2"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = "879,326.00"
$ws.Range("F6").Value = "151,331,741.00"
$ws.Range("G6").Value = "0.58"
$ws.Range("H6").Value = 603
$ws.Range("I6").Value = 14
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "4.88"
$ws.Range("L6").Value = "7.14"
$ws.Range("A7").Value = "Azizah Rahmawati"
$ws.Range("B7").Value = "Hansyah_S2l"
$ws.Range("C7").Value = "S2"
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = "6,511,624.00"
$ws.Range("F7").Value = "175,179,262.00"
$ws.Range("G7").Value = "3.72"
$ws.Range("H7").Value = 427
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = "0.00"
$ws.Range("L7").Value = "0.00"
$ws.Range("A8").Value = "Fadilah Damayanti"
$ws.Range("B8").Value = "Hansyah_S2l"
$ws.Range("C8").Value = "S2"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = "292,435.00"
$ws.Range("F8").Value = "179,487,985.00"
$ws.Range("G8").Value = "0.16"
$ws.Range("H8").Value = 405
$ws.Range("I8").Value = 14
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = "0.00"
$ws.Range("L8").Value = "0.00"
$ws.Range("A9").Value = "Annisa Putri Restu"
$ws.Range("B9").Value = "Hansyah_S2l"
$ws.Range("C9").Value = "S2"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "1,090,717.00"
$ws.Range("F9").Value = "186,099,111.00"
$ws.Range("G9").Value = "0.59"
$ws.Range("H9").Value = 1.42
$ws.Range("I9").Value = 14
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = "4.59"
$ws.Range("L9").Value = "7.14"
$ws.Range("A10").Value = "Riska Nurlita"
$ws.Range("B10").Value = "Hansyah_S2l"
$ws.Range("C10").Value = "S2"
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = "2,306,580.00"
$ws.Range("F10").Value = "188,817,452.00"
$ws.Range("G10").Value = "1.22"
$ws.Range("H10").Value = 552
$ws.Range("I10").Value = 14
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = "0.00"
$ws.Range("L10").Value = "0.00"
$ws.Range("A11").Value = "Debora Retima Sihombing"
$ws.Range("B11").Value = "Hansyah_S2l"
$ws.Range("C11").Value = "S2"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = "300,000.00"
$ws.Range("F11").Value = "166,650,859.00"
$ws.Range("G11").Value = "0.18"
$ws.Range("H11").Value = 825
$ws.Range("I11").Value = 14
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = "3.15"
$ws.Range("L11").Value = "0.00"
$ws.Range("A12").Value = "Erlangga Hutama"
$ws.Range("B12").Value = "Hansyah_S2l"
$ws.Range("C12").Value = "S2"
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = "587,668.00"
$ws.Range("F12").Value = "113,166,920.00"
$ws.Range("G12").Value = "0.52"
$ws.Range("H12").Value = 845
$ws.Range("I12").Value = 14
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = "0.00"
$ws.Range("L12").Value = "0.00"
$ws.Range("A13").Value = "Erick Ervan Dewanggga"
$ws.Range("B13").Value = "Hansyah_S2l"
$ws.Range("C13").Value = "S2"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = "0.00"
$ws.Range("F13").Value = "154,603,507.00"
$ws.Range("G13").Value = "0.00"
$ws.Range("H13").Value = 417
$ws.Range("I13").Value = 15
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = "0.00"
$ws.Range("L13").Value = "0.00"
$ws.Range("A14").Value = "Romli"
$ws.Range("B14").Value = "Hansyah_S2l"
$ws.Range("C14").Value = "S2"
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "0.00"
$ws.Range("F14").Value = "163,146,299.00"
$ws.Range("G14").Value = "0.00"
$ws.Range("H14").Value = 1.42
$ws.Range("I14").Value = 14
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = "0.00"
$ws.Range("L14").Value = "0.00"
$ws.Range("A15").Value = "Aldi Taufik"
$ws.Range("B15").Value = "Hansyah_S2l"
$ws.Range("C15").Value = "S2"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "410,638.00"
$ws.Range("F15").Value = "153,419,149.00"
$ws.Range("G15").Value = "0.27"
$ws.Range("H15").Value = 1.268
$ws.Range("I15").Value = 15
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = "4.15"
$ws.Range("L15").Value = "6.67"
$ws.Range("A16").Value = "Nur Halim"
$ws.Range("B16").Value = "Hansyah_S2l"
$ws.Range("C16").Value = "S2"
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = "1,729,613.00"
$ws.Range("F16").Value = "142,201,981.00"
$ws.Range("G16").Value = "1.22"
$ws.Range("H16").Value = 629
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = "5.16"
$ws.Range("L16").Value = "0.00"
$ws.Range("A17").Value = "Axl Wicaksono"
$ws.Range("B17").Value = "Hansyah_S2l"
$ws.Range("C17").Value = "S2"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = "0.00"
$ws.Range("F17").Value = "123,481,323.00"
$ws.Range("G17").Value = "0.00"
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 14
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = "0.00"
$ws.Range("L17").Value = "0.00"
$ws.Range("A18").Value = "Nuraini"
$ws.Range("B18").Value = "Hansyah_S2l"
$ws.Range("C18").Value = "S2"
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = "0.00"
$ws.Range("F18").Value = "102,605,709.00"
$ws.Range("G18").Value = "0.00"
$ws.Range("H18").Value = 328
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = "0.00"
$ws.Range("L18").Value = "0.00"

# Reset style back to Normal (removes the temporary @ text format
# we applied above) so the saved file has no extraneous cell styles,
# matching the original single-default-style workbook.
$ws.Range("A1:L18").Style = "Normal"

# Rename the sheet per the commit.
$ws.Name = "repayment_20250916_20250916 (1)"